$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "K"
$ws.Range("B4").Value = "L"
$ws.Range("C4").Value = "M"
$ws.Range("D4").Value = "N"
$ws.Range("E4").Value = "O"

$ws.Range("E4").Select()
